$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with the new daily records (rows 375-385),
# reusing the formatting of the last existing row (374).
$ws.Range("A374:D374").Copy() | Out-Null
$ws.Range("A375:D385").PasteSpecial(-4122) | Out-Null

$data = @(
    @(375, 44449, 1, 6, 173.2601790355183),
    @(376, 44450, 0, 6, 173.2601790355183),
    @(377, 44451, 3, 9, 259.8902685532775),
    @(378, 44452, 0, 9, 259.8902685532775),
    @(379, 44453, 0, 6, 173.2601790355183),
    @(380, 44454, 0, 6, 173.2601790355183),
    @(381, 44455, 0, 4, 115.5067860236789),
    @(382, 44456, 0, 3, 86.63008951775916),
    @(383, 44457, 0, 3, 86.63008951775916),
    @(384, 44458, 0, 0, 0),
    @(385, 44459, 0, 0, 0)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}
